# Re-sort the group_id / number-of-issues table into ascending numeric
# group order (group1, group2, ... group14) instead of the previous
# "sorted as text" order (group1, group10, group11, group12, ...), and
# leave the selection on the cell the user ended up clicking (R7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$groups = @("group1","group2","group3","group4","group5","group6","group7","group8","group9","group10","group11","group12","group13","group14")
$values = @(37,32,58,98,56,86,80,89,51,81,31,56,53,96)

for ($i = 0; $i -lt $groups.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $groups[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$null = $ws.Range("R7").Select()
